# Apply the "added font module and counter with font module" edit:
#  1. Row 13 (Camera movement) Total hours: "3 hours" -> "2 hours"
#  2. Row 24 task text: "Updated Camera movement and reset" -> "Updated Camera movement and reset camera"
#  3. New row 29: Death animation | Juan Hernández | 2 hours | 4 hours
#  4. New row 30: Round points | Martí Torres | 30 minutes | 30 minutes
#  5. Selection moved to I18:I19 (active cell I18)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update total hours for "Camera movement" row (row 13)
$ws.Range("E13").Value = "2 hours"

# 2. Rename task text in row 24
$ws.Range("B24").Value = "Updated Camera movement and reset camera"

# 3. Add new row 29
$ws.Range("B29").Value = "Death animation"
$ws.Range("C29").Value = "Juan Hernández"
$ws.Range("D29").Value = "2 hours"
$ws.Range("E29").Value = "4 hours"

# 4. Add new row 30
$ws.Range("B30").Value = "Round points"
$ws.Range("C30").Value = "Martí Torres"
$ws.Range("D30").Value = "30 minutes"
$ws.Range("E30").Value = "30 minutes"

# 5. Move the active selection, matching the saved view state
$ws.Range("I18:I19").Select()
